$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 63

# Columns A-D hold text (Date/Time/Weekday/Week formatted as strings, not
# native Excel date/time/number types) in the source data, matching every
# prior row in the sheet. Force Text number format before assigning so the
# COM layer doesn't auto-coerce date/time-looking or numeric-looking text
# into real date serials / numbers, then clear the format override so the
# new row keeps the same (unstyled) look as the rest of the table.
$textCols = 1,2,3,4
foreach ($c in $textCols) {
    $ws.Cells.Item($row, $c).NumberFormat = "@"
}

$ws.Cells.Item($row, 1).Value = "2023-06-20"
$ws.Cells.Item($row, 2).Value = "19:35:28"
$ws.Cells.Item($row, 3).Value = "Tuesday"
$ws.Cells.Item($row, 4).Value = "25"

foreach ($c in $textCols) {
    $ws.Cells.Item($row, $c).ClearFormats()
}

# Columns E-T are plain numbers.
$ws.Cells.Item($row, 5).Value = 122171
$ws.Cells.Item($row, 6).Value = 133686
$ws.Cells.Item($row, 7).Value = 162480
$ws.Cells.Item($row, 8).Value = 133267
$ws.Cells.Item($row, 9).Value = 177328
$ws.Cells.Item($row, 10).Value = 114631
$ws.Cells.Item($row, 11).Value = 201682
$ws.Cells.Item($row, 12).Value = 225461
$ws.Cells.Item($row, 13).Value = 175471
$ws.Cells.Item($row, 14).Value = 103847
$ws.Cells.Item($row, 15).Value = 39278
$ws.Cells.Item($row, 16).Value = 33884
$ws.Cells.Item($row, 17).Value = 51954
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36331
$ws.Cells.Item($row, 20).Value = -1
